$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo + merge mutable-residue rows into one row per block ---

# Column block A/B (SCMF/TRBP/K* rows use columns A & B)
$ws.Range("A1").Value = "Mutable Residues"
$ws.Range("A2").Value = "ARG-419 THR-702"
$ws.Range("A3").Value = ""

# Column block D/E
$ws.Range("D1").Value = "Mutable Residues"
$ws.Range("D2").Value = "ARG-419 MET-711"
$ws.Range("D3").Value = ""

# Column block G/H
$ws.Range("G1").Value = "Mutable Residues"
$ws.Range("G2").Value = "LYS-421 ARG-645"
$ws.Range("G3").Value = ""

# --- New 4th data block in columns J/K (header only, no values) ---
$ws.Range("J1").Value = "Mutable Residues"

$ws.Range("A5").Value = "Algorithm"
$ws.Range("D5").Value = "Algorithm"
$ws.Range("G5").Value = "Algorithm"
$ws.Range("J5").Value = "Algorithm"
$ws.Range("K5").Value = "LogZ"

$ws.Range("J6").Value = "SCMF"
$ws.Range("J7").Value = "TRBP"
$ws.Range("J8").Value = "K* (minimized)"
$ws.Range("J9").Value = "K* (pairwise lb)"

# --- Updated numeric values ---
$ws.Range("H6").Value = 28.38
$ws.Range("H7").Value = 84.4

# --- Column width for new column J (target stored width 15.6640625 chars;
#     engine quantizes ColumnWidth to ~1/6-character steps, so 14.83 is the
#     closest achievable input) ---
$ws.Range("J1").ColumnWidth = 14.83

# --- Selection / view state ---
$ws.Range("K5").Select()
